# Weekly update: a new weekly price record for "Acelga" (Región Metropolitana,
# Vega Central Mapocho de Santiago) is inserted at row 333, pushing the
# existing rows 333:448 down to 334:449 (dimension grows from R448 to R449).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 333, shifting the rest of the table down.
$ws.Rows.Item(333).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A333").Value = 9
$ws.Range("B333").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C333").Value = "Metropolitana"
$ws.Range("D333").Value = 44559
$ws.Range("E333").Value = 13
$ws.Range("F333").Value = 100112009
$ws.Range("G333").Value = "Acelga"
$ws.Range("H333").Value = "Sin especificar"
$ws.Range("I333").Value = "Primera"
$ws.Range("J333").Value = 52
$ws.Range("K333").Value = 14000
$ws.Range("L333").Value = 15000
$ws.Range("M333").Value = 14500
$ws.Range("N333").Value = "`$/docena de atados"
$ws.Range("O333").Value = "Región Metropolitana"
$ws.Range("P333").Value = 4833
$ws.Range("Q333").Value = 3
$ws.Range("R333").Value = "Hortaliza"

# Keep the same date number format used by the other rows in column D.
$ws.Range("D333").NumberFormat = $ws.Range("D334").NumberFormat
